# EIA Table 2.5.A — roll the report forward one month: October 2016 -> November 2016.
# Adds a "November" monthly data row under the "Year 2016" block (row 53, pushing the
# "Year to Date" / "Rolling 12 Months" summary rows down by one), refreshes the
# Year-to-Date and Rolling-12-Months totals, and updates the title / rolling-window
# caption text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Make room for the new "November" row at row 53 (currently the merged
#    "Year to Date" label row). Unmerge it first so Insert() cleanly shifts
#    every cell below down by one row (merged ranges further down shift on
#    their own).
# ---------------------------------------------------------------------------
$ws.Range("A53:F53").UnMerge()
$ws.Rows("53:53").Insert(-4121)   # xlShiftDown
$ws.Range("A54:F54").Merge()      # restore the "Year to Date" label merge at its new row

# Give the freshly inserted row 53 the same formatting as a normal monthly
# data row (copy the look of row 52, "October") instead of the blank format
# Insert() minted for it.
$ws.Range("A52:F52").Copy()
$ws.Range("A53:F53").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. Fill in the new "November" monthly row.
# ---------------------------------------------------------------------------
$ws.Range("A53").Value2 = "November"
$ws.Range("B53").Value2 = 26480
$ws.Range("C53").Value2 = 2417
$ws.Range("D53").Value2 = 21537
$ws.Range("E53").Value2 = 2091
$ws.Range("F53").Value2 = 436

# ---------------------------------------------------------------------------
# 3. Refresh the "Year to Date" totals (now rows 55-57: 2014, 2015, 2016).
# ---------------------------------------------------------------------------
$ws.Range("B55").Value2 = 263531
$ws.Range("C55").Value2 = 23757
$ws.Range("D55").Value2 = 210555
$ws.Range("E55").Value2 = 24934
$ws.Range("F55").Value2 = 4285

$ws.Range("B56").Value2 = 257009
$ws.Range("C56").Value2 = 23092
$ws.Range("D56").Value2 = 206578
$ws.Range("E56").Value2 = 23135
$ws.Range("F56").Value2 = 4204

$ws.Range("B57").Value2 = 301417
$ws.Range("C57").Value2 = 27984
$ws.Range("D57").Value2 = 245471
$ws.Range("E57").Value2 = 23173
$ws.Range("F57").Value2 = 4789

# ---------------------------------------------------------------------------
# 4. Refresh the "Rolling 12 Months" totals (now rows 59-60: 2015, 2016) and
#    its caption (row 58).
# ---------------------------------------------------------------------------
$ws.Range("A58").Value2 = "Rolling 12 Months Ending in November"

$ws.Range("B59").Value2 = 279461
$ws.Range("C59").Value2 = 25154
$ws.Range("D59").Value2 = 224470
$ws.Range("E59").Value2 = 25240
$ws.Range("F59").Value2 = 4597

$ws.Range("B60").Value2 = 326937
$ws.Range("C60").Value2 = 30150
$ws.Range("D60").Value2 = 266273
$ws.Range("E60").Value2 = 25288
$ws.Range("F60").Value2 = 5227

# ---------------------------------------------------------------------------
# 5. Update the report title/caption to reference November instead of October.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value2 = "by Sector, 2006-November 2016 (Million Cubic Feet)"
